$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the color codes in column B for the rows whose values changed.
$ws.Range("B3").Value  = "#2A4845"
$ws.Range("B4").Value  = "#E94E24"
$ws.Range("B5").Value  = "#579CB4"
$ws.Range("B6").Value  = "#FB8F71"
$ws.Range("B7").Value  = "#ED6D63"
$ws.Range("B9").Value  = "#41B496"
$ws.Range("B10").Value = "#75C7C7"
$ws.Range("B11").Value = "#FFBA6E"
$ws.Range("B13").Value = "#93B9B9"
$ws.Range("B15").Value = "#F4A28C"
$ws.Range("B16").Value = "#F6B3AF"
$ws.Range("B18").Value = "#8AD3FB"

# Move the active selection to B19 (matches the recorded UI state in the edit).
$ws.Range("B19").Select()
